$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update environment / credentials / account values (row 2) ---
$ws.Range("A2").Value = "ssurgwsoadev4-oci.opc.oracleoutsourcing.com"
$ws.Range("B2").Value = "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/pc/PolicyCenter.do"
$ws.Range("D2").Value = "gw"
$ws.Range("E2").Value = 9498924883
$ws.Range("J2").Value = "No"

# --- Remove the hyperlink that used to live on B2 (old environment URL) ---
$ws.Range("B2").Hyperlinks.Delete()

# --- Clear the FechaInicio value (Q2), keeping its number format/style ---
$ws.Range("Q2").ClearContents()

# --- Replace the external-workbook formulas with static values ---
$ws.Range("W2").Value = "MMM111"
$ws.Range("X2").Value = "MASDAS12312"
$ws.Range("Y2").Value = "ASDAKE1232"

# --- Break the external link now that nothing references it anymore ---
$links = $wb.LinkSources()
if ($links) {
    foreach ($link in $links) {
        $wb.BreakLink($link, 1)
    }
}

# --- Update the active selection to match the new state ---
$ws.Range("W3").Select()
